$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 836.6667
$ws.Range("I20").Value = 836.6667
$ws.Range("K20").Value = 836.6667
$ws.Range("M20").Value = -606.6667
$ws.Range("H35").Value = 836.6667
$ws.Range("I35").Value = 836.6667
$ws.Range("K35").Value = 836.6667
$ws.Range("M35").Value = -457.6667
$ws.Range("H42").Value = 120.27273
$ws.Range("I42").Value = 66.333336
$ws.Range("J42").Value = 185
$ws.Range("K42").Value = 199.000008
$ws.Range("L42").Value = 555
$ws.Range("M42").Value = 30.99999199999999
$ws.Range("N42").Value = -1015
$ws.Range("H105").Value = 60000
$ws.Range("I105").Value = 60000
$ws.Range("K105").Value = 60000
$ws.Range("M105").Value = -56506
$ws.Range("H107").Value = 425.33334
$ws.Range("I107").Value = 425.33334
$ws.Range("K107").Value = 425.33334
$ws.Range("M107").Value = 1494.66666
$ws.Range("I113").Value = 4000
$ws.Range("K113").Value = 4000
$ws.Range("M113").Value = -746
$ws.Range("H116").Value = 6540.1665
$ws.Range("I116").Value = 6631
$ws.Range("K116").Value = 6631
$ws.Range("M116").Value = -3189
$ws.Range("H137").Value = 1861.3182
$ws.Range("I137").Value = 1750.6471
$ws.Range("K137").Value = 5251.9413
$ws.Range("M137").Value = -2701.9413
$ws.Range("H138").Value = 3750.1177
$ws.Range("I138").Value = 941.6923
$ws.Range("K138").Value = 2825.0769
$ws.Range("M138").Value = 2314.9231

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2952.5833
$ws.Range("I2").Value = 1979.25
$ws.Range("K2").Value = 1979.25
$ws.Range("M2").Value = -1866.25
$ws.Range("H32").Value = 4404.5884
$ws.Range("I32").Value = 3351.7026
$ws.Range("J32").Value = 7187.2144
$ws.Range("K32").Value = 3351.7026
$ws.Range("L32").Value = 7187.2144
$ws.Range("M32").Value = -3064.7026
$ws.Range("N32").Value = -7761.2144
$ws.Range("H61").Value = 1766.32
$ws.Range("J61").Value = 2893.75
$ws.Range("L61").Value = 2893.75
$ws.Range("N61").Value = -3317.75
$ws.Range("H110").Value = 4591
$ws.Range("I110").Value = 4767.778
$ws.Range("K110").Value = 4767.778
$ws.Range("M110").Value = -2722.778
$ws.Range("H116").Value = 2952.5833
$ws.Range("I116").Value = 1979.25
$ws.Range("K116").Value = 1979.25
$ws.Range("M116").Value = 314.75
$ws.Range("H122").Value = 373100.62
$ws.Range("J122").Value = 2777.4
$ws.Range("L122").Value = 8332.200000000001
$ws.Range("N122").Value = -13232.2
$ws.Range("H136").Value = 1766.32
$ws.Range("J136").Value = 2893.75
$ws.Range("L136").Value = 8681.25
$ws.Range("N136").Value = -13781.25
$ws.Range("H140").Value = 84831.664
$ws.Range("J140").Value = 88798
$ws.Range("L140").Value = 88798
$ws.Range("N140").Value = -99158

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2952.5833
$ws.Range("I3").Value = 1979.25
$ws.Range("K3").Value = 1979.25
$ws.Range("M3").Value = -1865.25
$ws.Range("H55").Value = 34999
$ws.Range("J55").Value = 34999
$ws.Range("L55").Value = 34999
$ws.Range("N55").Value = -35545
$ws.Range("H99").Value = 34801.8
$ws.Range("I99").Value = 44819.26
$ws.Range("J99").Value = 1887.2858
$ws.Range("K99").Value = 44819.26
$ws.Range("L99").Value = 1887.2858
$ws.Range("M99").Value = -43321.26
$ws.Range("N99").Value = -4883.2858
$ws.Range("H105").Value = 2759.2
$ws.Range("I105").Value = 2599.1428
$ws.Range("K105").Value = 2599.1428
$ws.Range("M105").Value = -852.1428000000001
$ws.Range("H134").Value = 2014.3
$ws.Range("I134").Value = 1701.826
$ws.Range("K134").Value = 5105.478
$ws.Range("M134").Value = -2570.478

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3730.932
$ws.Range("I31").Value = 3139.1667
$ws.Range("J31").Value = 4140.615
$ws.Range("K31").Value = 3139.1667
$ws.Range("L31").Value = 4140.615
$ws.Range("M31").Value = -2844.1667
$ws.Range("N31").Value = -4730.615
$ws.Range("H34").Value = 3730.932
$ws.Range("I34").Value = 3139.1667
$ws.Range("J34").Value = 4140.615
$ws.Range("K34").Value = 3139.1667
$ws.Range("L34").Value = 4140.615
$ws.Range("M34").Value = -2937.1667
$ws.Range("N34").Value = -4544.615
$ws.Range("H58").Value = 2794.5881
$ws.Range("I58").Value = 1268.35
$ws.Range("J58").Value = 4974.9287
$ws.Range("K58").Value = 1268.35
$ws.Range("L58").Value = 4974.9287
$ws.Range("M58").Value = -1065.35
$ws.Range("N58").Value = -5380.9287
$ws.Range("H62").Value = 43257.4
$ws.Range("I62").Value = 3264.5
$ws.Range("J62").Value = 103246.75
$ws.Range("K62").Value = 3264.5
$ws.Range("L62").Value = 103246.75
$ws.Range("M62").Value = -2640.5
$ws.Range("N62").Value = -104494.75
$ws.Range("H65").Value = 43257.4
$ws.Range("I65").Value = 3264.5
$ws.Range("J65").Value = 103246.75
$ws.Range("K65").Value = 16322.5
$ws.Range("L65").Value = 516233.75
$ws.Range("M65").Value = -13202.5
$ws.Range("N65").Value = -522473.75
$ws.Range("H99").Value = 14910.857
$ws.Range("I99").Value = 10865.75
$ws.Range("K99").Value = 10865.75
$ws.Range("M99").Value = -9367.75
$ws.Range("H107").Value = 1183.1765
$ws.Range("I107").Value = 881.7143
$ws.Range("J107").Value = 1394.2
$ws.Range("K107").Value = 881.7143
$ws.Range("L107").Value = 1394.2
$ws.Range("M107").Value = 1038.2857
$ws.Range("N107").Value = -5234.2
$ws.Range("H122").Value = 2427.3333
$ws.Range("I122").Value = 2477.7368
$ws.Range("J122").Value = 1948.5
$ws.Range("K122").Value = 7433.2104
$ws.Range("L122").Value = 5845.5
$ws.Range("M122").Value = -4983.2104
$ws.Range("N122").Value = -10745.5
$ws.Range("H126").Value = 14910.857
$ws.Range("I126").Value = 10865.75
$ws.Range("K126").Value = 32597.25
$ws.Range("M126").Value = -30127.25
$ws.Range("H134").Value = 2348.6445
$ws.Range("I134").Value = 1955.2122
$ws.Range("K134").Value = 5865.6366
$ws.Range("M134").Value = -3330.6366
$ws.Range("H136").Value = 2794.5881
$ws.Range("I136").Value = 1268.35
$ws.Range("J136").Value = 4974.9287
$ws.Range("K136").Value = 3805.05
$ws.Range("L136").Value = 14924.7861
$ws.Range("M136").Value = -1255.05
$ws.Range("N136").Value = -20024.7861

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H93").Value = 5198.4
$ws.Range("J93").Value = 4498
$ws.Range("L93").Value = 13494
$ws.Range("N93").Value = -17238
$ws.Range("H115").Value = 428
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
$ws.Range("H131").Value = 1494
$ws.Range("J131").Value = 1494
$ws.Range("L131").Value = 4482
$ws.Range("N131").Value = -14562

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 3750.75
$ws.Range("I13").Value = 5
$ws.Range("K13").Value = 5
$ws.Range("M13").Value = 134
$ws.Range("H97").Value = 1147.1875
$ws.Range("I97").Value = 1168.2142
$ws.Range("K97").Value = 1168.2142
$ws.Range("M97").Value = -672.2141999999999
$ws.Range("H122").Value = 79847.30499999999
$ws.Range("I122").Value = 2821.6
$ws.Range("K122").Value = 8464.799999999999
$ws.Range("M122").Value = -6014.799999999999
$ws.Range("H123").Value = 24325.143
$ws.Range("J123").Value = 24325.143
$ws.Range("L123").Value = 24325.143
$ws.Range("N123").Value = -29225.143

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2938.5833
$ws.Range("I40").Value = 3037.1
$ws.Range("J40").Value = 2446
$ws.Range("K40").Value = 3037.1
$ws.Range("L40").Value = 2446
$ws.Range("M40").Value = -2901.1
$ws.Range("N40").Value = -2718
$ws.Range("H100").Value = 2532.4
$ws.Range("I100").Value = 2887.3333
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 2887.3333
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -2346.3333
$ws.Range("N100").Value = -3082
$ws.Range("H132").Value = 4370.409
$ws.Range("I132").Value = 3122.7273
$ws.Range("K132").Value = 9368.1819
$ws.Range("M132").Value = -6838.1819

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 5125
$ws.Range("J4").Value = 5125
$ws.Range("L4").Value = 5125
$ws.Range("N4").Value = -5351
$ws.Range("H5").Value = 1510001
$ws.Range("I5").Value = 3000000
$ws.Range("J5").Value = 20002
$ws.Range("K5").Value = 3000000
$ws.Range("L5").Value = 20002
$ws.Range("M5").Value = -2999888
$ws.Range("N5").Value = -20226
$ws.Range("H58").Value = 10000
$ws.Range("I58").Value = 10000
$ws.Range("K58").Value = 10000
$ws.Range("M58").Value = -9692
$ws.Range("H74").Value = 60673.5
$ws.Range("J74").Value = 60673.5
$ws.Range("L74").Value = 60673.5
$ws.Range("N74").Value = -62545.5
$ws.Range("H77").Value = 60673.5
$ws.Range("J77").Value = 60673.5
$ws.Range("L77").Value = 182020.5
$ws.Range("N77").Value = -191380.5
$ws.Range("H126").Value = 2899.2727
$ws.Range("I126").Value = 2099.111
$ws.Range("J126").Value = 6500
$ws.Range("K126").Value = 6297.333
$ws.Range("L126").Value = 19500
$ws.Range("M126").Value = -3827.333
$ws.Range("N126").Value = -24440
$ws.Range("H132").Value = 38975.64
$ws.Range("I132").Value = 44027.047
$ws.Range("K132").Value = 132081.141
$ws.Range("M132").Value = -129551.141
